# Generate Report for Handoff
# Adds two new tracked files to the localization-status report:
#   7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md
#   c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Overview sheet (sheet1): columns A (File Name) B (zh-cn) C (de-de)
# D (Latest Handoff Date)
# ---------------------------------------------------------------------

# Row 4 - 7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md
$wsOverview.Range("A4").Value = "7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/7fe64e27ed125ae01df10791eb33567d14f92276/e2e/7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md", $null, $null, "7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md")
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
$wsOverview.Range("D4").Value = "2016-29-18 00:29:53"

# Row 5 - c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md
$wsOverview.Range("A5").Value = "c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c6cfb8e565e9ae01df10791eb33567d14f92276/e2e/c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md", $null, $null, "c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md")
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = "2016-29-18 00:29:53"

# ---------------------------------------------------------------------
# zh-cn sheet (sheet2): columns
# A Source File Name, B File Extension, C Status, D Latest Handoff File,
# E Latest Handoff Datetime, H Latest Handback DateTime, I Handoff Reason
# ---------------------------------------------------------------------

# Row 4 - 7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd
$wsZhCn.Range("A4").Value = "7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/7fe64e27ed125ae01df10791eb33567d14f92276/e2e/7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md", $null, $null, "7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md")
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/7fe64e27ed125ae01df10791eb33567d14f92276/e2e/7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md", $null, $null, ".md")
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.965bd93b8d7f0ff29d684f0a71d9bd2b1addfaa2.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/965bd93b8d7f0ff29d684f0a71d9bd2b1addfaa/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.965bd93b8d7f0ff29d684f0a71d9bd2b1addfaa2.zh-cn.xlf", $null, $null, "7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.965bd93b8d7f0ff29d684f0a71d9bd2b1addfaa2.zh-cn.xlf")
$wsZhCn.Range("E4").Value = "2016-03-18 00:29:50"
$wsZhCn.Range("E4").NumberFormat = $dateFormat
$wsZhCn.Range("H4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I4").Value = "Include"

# Row 5 - c6cfb8e5-65e9-48fd-bf95-3912fcae9701
$wsZhCn.Range("A5").Value = "c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c6cfb8e565e9ae01df10791eb33567d14f92276/e2e/c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md", $null, $null, "c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md")
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/c6cfb8e565e9ae01df10791eb33567d14f92276/e2e/c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md", $null, $null, ".md")
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "c6cfb8e5-65e9-48fd-bf95-3912fcae9701.1669b18482b43a3dca05dabb6b7313f28bf459fc.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1669b18482b43a3dca05dabb6b7313f28bf459f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c6cfb8e5-65e9-48fd-bf95-3912fcae9701.1669b18482b43a3dca05dabb6b7313f28bf459fc.zh-cn.xlf", $null, $null, "c6cfb8e5-65e9-48fd-bf95-3912fcae9701.1669b18482b43a3dca05dabb6b7313f28bf459fc.zh-cn.xlf")
$wsZhCn.Range("E5").Value = "2016-03-18 00:29:50"
$wsZhCn.Range("E5").NumberFormat = $dateFormat
$wsZhCn.Range("H5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I5").Value = "Include"

# ---------------------------------------------------------------------
# de-de sheet (sheet3): same column layout as zh-cn
# ---------------------------------------------------------------------

# Row 4 - 7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd
$wsDeDe.Range("A4").Value = "7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/7fe64e27ed125ae01df10791eb33567d14f92276/e2e/7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md", $null, $null, "7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md")
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/7fe64e27ed125ae01df10791eb33567d14f92276/e2e/7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.md", $null, $null, ".md")
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.965bd93b8d7f0ff29d684f0a71d9bd2b1addfaa2.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/965bd93b8d7f0ff29d684f0a71d9bd2b1addfaa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.965bd93b8d7f0ff29d684f0a71d9bd2b1addfaa2.de-de.xlf", $null, $null, "7fe64e27-08d4-4d3e-8271-ecf79fe1e9fd.965bd93b8d7f0ff29d684f0a71d9bd2b1addfaa2.de-de.xlf")
$wsDeDe.Range("E4").Value = "2016-03-18 00:29:53"
$wsDeDe.Range("E4").NumberFormat = $dateFormat
$wsDeDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I4").Value = "Include"

# Row 5 - c6cfb8e5-65e9-48fd-bf95-3912fcae9701
$wsDeDe.Range("A5").Value = "c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c6cfb8e565e9ae01df10791eb33567d14f92276/e2e/c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md", $null, $null, "c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md")
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/c6cfb8e565e9ae01df10791eb33567d14f92276/e2e/c6cfb8e5-65e9-48fd-bf95-3912fcae9701.md", $null, $null, ".md")
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "c6cfb8e5-65e9-48fd-bf95-3912fcae9701.1669b18482b43a3dca05dabb6b7313f28bf459fc.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1669b18482b43a3dca05dabb6b7313f28bf459f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c6cfb8e5-65e9-48fd-bf95-3912fcae9701.1669b18482b43a3dca05dabb6b7313f28bf459fc.de-de.xlf", $null, $null, "c6cfb8e5-65e9-48fd-bf95-3912fcae9701.1669b18482b43a3dca05dabb6b7313f28bf459fc.de-de.xlf")
$wsDeDe.Range("E5").Value = "2016-03-18 00:29:53"
$wsDeDe.Range("E5").NumberFormat = $dateFormat
$wsDeDe.Range("H5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I5").Value = "Include"
